# Update the cryptos price/volume table with the latest scraped values.
# Column D ("Price") values are text (dotted-thousands formatted numbers,
# e.g. "58.035.03"), so each is written with a leading apostrophe to force
# text entry and stop Excel from re-interpreting it as a number (which
# would silently drop things like trailing zeros, e.g. 517.70 -> 517.7).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "'58.035.03"
$ws.Range("E2").Value2 = "  -1.37%  "
$ws.Range("D3").Value2 = "'2.464.26"
$ws.Range("E3").Value2 = "  -1.32%  "
$ws.Range("D4").Value2 = "'0.999"
$ws.Range("E4").Value2 = "  -0.41%  "
$ws.Range("D5").Value2 = "'517.70"
$ws.Range("E5").Value2 = "  -3.38%  "
$ws.Range("D6").Value2 = "'131.17"
$ws.Range("E6").Value2 = "  -3.77%  "
$ws.Range("D7").Value2 = "'0.999"
$ws.Range("E7").Value2 = "  -0.02%  "
$ws.Range("E8").Value2 = "  -1.79%  "
$ws.Range("D9").Value2 = "'2.466.67"
$ws.Range("E9").Value2 = "  -2.33%  "
$ws.Range("E10").Value2 = "  -2.22%  "
$ws.Range("D11").Value2 = "'0.157"
$ws.Range("E11").Value2 = "  -0.22%  "
$ws.Range("E12").Value2 = "  -0.15%  "
$ws.Range("D13").Value2 = "'0.339"
$ws.Range("E13").Value2 = "  -2.42%  "
$ws.Range("D14").Value2 = "'2.895.04"
$ws.Range("E14").Value2 = "  -2.53%  "
$ws.Range("D15").Value2 = "'57.951.31"
$ws.Range("E15").Value2 = "  -1.68%  "
$ws.Range("D16").Value2 = "'22.27"
$ws.Range("E16").Value2 = "  -3.10%  "
$ws.Range("E17").Value2 = "  -2.25%  "
$ws.Range("D18").Value2 = "'2.442.74"
$ws.Range("E18").Value2 = "  -2.98%  "
$ws.Range("D19").Value2 = "'10.70"
$ws.Range("E19").Value2 = "  -3.72%  "
$ws.Range("D20").Value2 = "'319.90"
$ws.Range("E20").Value2 = "  -1.15%  "
$ws.Range("D21").Value2 = "'4.15"
$ws.Range("E21").Value2 = "  -2.57%  "
$ws.Range("E22").Value2 = "  +0.02%  "
$ws.Range("D23").Value2 = "'5.71"
$ws.Range("E23").Value2 = "  -3.74%  "
$ws.Range("D24").Value2 = "'64.09"
$ws.Range("E24").Value2 = "  -1.69%  "
$ws.Range("D25").Value2 = "'0.407"
$ws.Range("E25").Value2 = "  -2.93%  "
$ws.Range("E26").Value2 = "  +0.25%  "
$ws.Range("E27").Value2 = "  -3.55%  "
$ws.Range("D28").Value2 = "'7.30"
$ws.Range("E28").Value2 = "  -2.52%  "
$ws.Range("D29").Value2 = "'0.0₃0738"
$ws.Range("E29").Value2 = "  -4.03%  "
$ws.Range("D30").Value2 = "'165.69"
$ws.Range("E30").Value2 = "  -3.53%  "
$ws.Range("E31").Value2 = "  -4.28%  "
$ws.Range("E32").Value2 = "  -6.22%  "
$ws.Range("D33").Value2 = "'1.16"
$ws.Range("E33").Value2 = "  -0.77%  "
$ws.Range("D34").Value2 = "'0.998"
$ws.Range("D35").Value2 = "'0.998"
$ws.Range("E35").Value2 = "  +0.33%  "
$ws.Range("E36").Value2 = "  -2.12%  "
$ws.Range("D37").Value2 = "'1.29"
$ws.Range("E37").Value2 = "  -8.19%  "
$ws.Range("E38").Value2 = "  -3.35%  "
$ws.Range("D39").Value2 = "'1.47"
$ws.Range("E39").Value2 = "  -4.64%  "
$ws.Range("D40").Value2 = "'0.787"
$ws.Range("E40").Value2 = "  -2.79%  "
$ws.Range("B41").Value2 = "RenderToken"
$ws.Range("C41").Value2 = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").Value2 = "'5.11"
$ws.Range("E41").Value2 = "  -1.40%  "
$ws.Range("B42").Value2 = "Filecoin"
$ws.Range("C42").Value2 = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").Value2 = "'3.44"
$ws.Range("E42").Value2 = "  -4.08%  "
$ws.Range("B43").Value2 = "Bittensor"
$ws.Range("C43").Value2 = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D43").Value2 = "'270.30"
$ws.Range("E43").Value2 = "  -5.45%  "
$ws.Range("D44").Value2 = "'0.591"
$ws.Range("E44").Value2 = "  -2.96%  "
$ws.Range("D45").Value2 = "'125.08"
$ws.Range("E45").Value2 = "  -4.64%  "
$ws.Range("D46").Value2 = "'0.0905"
$ws.Range("E46").Value2 = "  -1.82%  "
$ws.Range("E47").Value2 = "  -4.01%  "
$ws.Range("E48").Value2 = "  -4.97%  "
$ws.Range("E49").Value2 = "  -3.35%  "
$ws.Range("D50").Value2 = "'1.718.10"
$ws.Range("E50").Value2 = "  -2.27%  "
$ws.Range("D51").Value2 = "'0.967"
$ws.Range("E51").Value2 = "  -2.45%  "
